# Chiffres COVID-19 Valais.xlsx - daily data update
# Updates the "Nb nouveaux cas positifs" (col C), a couple of
# hospital/ICU snapshot columns (E/F/G) and the death-breakdown columns
# (L/M) for the days around 2021-01-04..01-07 (rows 314-317), and fills
# in the newly-reported day (row 317, 2021-01-07) that was previously
# blank. Columns B, H, J and K are running-total formulas and recompute
# automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 308 (2020-12-29): new positive cases 137 -> 138
$ws.Range("C308").Value = 138

# Row 309 (2020-12-30): new positive cases 143 -> 144
$ws.Range("C309").Value = 144

# Row 311 (2021-01-01): new positive cases 72 -> 73
$ws.Range("C311").Value = 73

# Row 312 (2021-01-02): new positive cases 149 -> 148
$ws.Range("C312").Value = 148

# Row 314 (2021-01-04): new positive cases 236 -> 237
$ws.Range("C314").Value = 237

# Row 315 (2021-01-05): new positive cases 124 -> 170, extra-hospital deaths 1 -> 2
$ws.Range("C315").Value = 170
$ws.Range("M315").Value = 2

# Row 316 (2021-01-06): new positive cases 22 -> 115, extra-hospital deaths 0 -> 2
$ws.Range("C316").Value = 115
$ws.Range("M316").Value = 2

# Row 317 (2021-01-07): previously-blank day now fully reported
$ws.Range("C317").Value = 31
$ws.Range("E317").Value = 10
$ws.Range("F317").Value = 8
$ws.Range("G317").Value = 83
$ws.Range("L317").Value = 0
$ws.Range("M317").Value = 0

# Move the frozen-pane view / active selection down to the newly edited rows
$ws.Range("S35").Select()
